$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.06649561763787659
$ws.Cells.Item(2, 4).Value = 0.1603103226534159
$ws.Cells.Item(2, 5).Value = 0.07830362079268482
$ws.Cells.Item(2, 6).Value = 2.761593524215925
$ws.Cells.Item(2, 7).Value = 0.002579208089904749
$ws.Cells.Item(2, 9).Value = 2.524608240905252
$ws.Cells.Item(2, 11).Value = 2.083135759181971
$ws.Cells.Item(2, 12).Value = 0.1324669927932725
$ws.Cells.Item(2, 13).Value = 0.5559387455636937
$ws.Cells.Item(3, 3).Value = 0.06624426197548416
$ws.Cells.Item(3, 4).Value = 0.16149832773025
$ws.Cells.Item(3, 5).Value = 0.07823285078791997
$ws.Cells.Item(3, 6).Value = 2.731702127650692
$ws.Cells.Item(3, 7).Value = 0.002584813209828425
$ws.Cells.Item(3, 9).Value = 2.504424526509567
$ws.Cells.Item(3, 11).Value = 1.96735084434971
$ws.Cells.Item(3, 12).Value = 0.1325066693301373
$ws.Cells.Item(3, 13).Value = 0.5348443598449464
$ws.Cells.Item(4, 3).Value = 0.06611379606925993
$ws.Cells.Item(4, 4).Value = 0.1622712438345104
$ws.Cells.Item(4, 5).Value = 0.07821993004529304
$ws.Cells.Item(4, 6).Value = 2.714908764533206
$ws.Cells.Item(4, 7).Value = 0.002588433836526529
$ws.Cells.Item(4, 9).Value = 2.493390041777033
$ws.Cells.Item(4, 11).Value = 1.897504434658998
$ws.Cells.Item(4, 12).Value = 0.1325844558708127
$ws.Cells.Item(4, 13).Value = 0.5222258120514596
$ws.Cells.Item(5, 3).Value = 0.06606666338279865
$ws.Cells.Item(5, 4).Value = 0.162597110240327
$ws.Cells.Item(5, 5).Value = 0.07822234305127829
$ws.Cells.Item(5, 6).Value = 2.708455430729714
$ws.Cells.Item(5, 7).Value = 0.00258995445609898
$ws.Cells.Item(5, 9).Value = 2.489233263750421
$ws.Cells.Item(5, 11).Value = 1.869353997739353
$ws.Cells.Item(5, 12).Value = 0.1326295632877788
$ws.Cells.Item(5, 13).Value = 0.5171673721863357
$ws.Cells.Item(6, 3).Value = 0.06605920236756901
$ws.Cells.Item(6, 4).Value = 0.1626518772109478
$ws.Cells.Item(6, 5).Value = 0.07822320750869238
$ws.Cells.Item(6, 6).Value = 2.707407356550817
$ws.Cells.Item(6, 7).Value = 0.002590209687444234
$ws.Cells.Item(6, 9).Value = 2.488563510888326
$ws.Cells.Item(6, 11).Value = 1.864698494550794
$ws.Cells.Item(6, 12).Value = 0.1326378624621611
$ws.Cells.Item(6, 13).Value = 0.5163324748472178
$ws.Cells.Item(7, 3).Value = 0.06611313595028889
$ws.Cells.Item(7, 4).Value = 0.1622755945106498
$ws.Cells.Item(7, 5).Value = 0.07821993149731377
$ws.Cells.Item(7, 6).Value = 2.714820156054415
$ws.Cells.Item(7, 7).Value = 0.002588454161068619
$ws.Cells.Item(7, 9).Value = 2.49333260820984
$ws.Cells.Item(7, 11).Value = 1.897123523541694
$ws.Cells.Item(7, 12).Value = 0.1325850099466237
$ws.Cells.Item(7, 13).Value = 0.5221572533562906
$ws.Cells.Item(8, 3).Value = 0.06640401505599414
$ws.Cells.Item(8, 4).Value = 0.1607109019676543
$ws.Cells.Item(8, 5).Value = 0.07827288245864672
$ws.Cells.Item(8, 6).Value = 2.750961837763256
$ws.Cells.Item(8, 7).Value = 0.002581103667098179
$ws.Cells.Item(8, 9).Value = 2.517365914050288
$ws.Cells.Item(8, 11).Value = 2.042953990591627
$ws.Cells.Item(8, 12).Value = 0.1324695692967062
$ws.Cells.Item(8, 13).Value = 0.5485960909962984
$ws.Cells.Item(9, 3).Value = 0.06716253052544374
$ws.Cells.Item(9, 4).Value = 0.1579890664498329
$ws.Cells.Item(9, 5).Value = 0.07861904982067713
$ws.Cells.Item(9, 6).Value = 2.834319825688226
$ws.Cells.Item(9, 7).Value = 0.002568102885342613
$ws.Cells.Item(9, 9).Value = 2.575356072527441
$ws.Cells.Item(9, 11).Value = 2.338868504906145
$ws.Cells.Item(9, 12).Value = 0.1326683669336504
$ws.Cells.Item(9, 13).Value = 0.6030988710625778
$ws.Cells.Item(10, 3).Value = 0.06783293315621108
$ws.Cells.Item(10, 4).Value = 0.1562026704976525
$ws.Cells.Item(10, 5).Value = 0.0790213568048479
$ws.Cells.Item(10, 6).Value = 2.903329865980538
$ws.Cells.Item(10, 7).Value = 0.002559402754675988
$ws.Cells.Item(10, 9).Value = 2.624704577132107
$ws.Cells.Item(10, 11).Value = 2.562442644498503
$ws.Cells.Item(10, 12).Value = 0.1330755724138726
$ws.Cells.Item(10, 13).Value = 0.6447812977723402
$ws.Cells.Item(11, 3).Value = 0.06816217095049382
$ws.Cells.Item(11, 4).Value = 0.1554367704833446
$ws.Cells.Item(11, 5).Value = 0.07923657771191372
$ws.Cells.Item(11, 6).Value = 2.936445353842458
$ws.Cells.Item(11, 7).Value = 0.0025556275691872
$ws.Cells.Item(11, 9).Value = 2.648645685294809
$ws.Cells.Item(11, 11).Value = 2.665516571846467
$ws.Cells.Item(11, 12).Value = 0.1333179666745892
$ws.Cells.Item(11, 13).Value = 0.6641044933817568
$ws.Cells.Item(12, 3).Value = 0.06829030719728735
$ws.Cells.Item(12, 4).Value = 0.1551535055901034
$ws.Cells.Item(12, 5).Value = 0.07932271181643102
$ws.Cells.Item(12, 6).Value = 2.949235661564614
$ws.Cells.Item(12, 7).Value = 0.002554224088293706
$ws.Cells.Item(12, 9).Value = 2.657928213287576
$ws.Cells.Item(12, 11).Value = 2.704746444338014
$ws.Cells.Item(12, 12).Value = 0.133418007394134
$ws.Cells.Item(12, 13).Value = 0.6714739937772833
$ws.Cells.Item(13, 3).Value = 0.06826255728337571
$ws.Cells.Item(13, 4).Value = 0.1552142101855303
$ws.Cells.Item(13, 5).Value = 0.07930395511908372
$ws.Cells.Item(13, 6).Value = 2.946469872527587
$ws.Cells.Item(13, 7).Value = 0.002554525194561785
$ws.Cells.Item(13, 9).Value = 2.655919392097658
$ws.Cells.Item(13, 11).Value = 2.696288763449502
$ws.Cells.Item(13, 12).Value = 0.1333960943334702
$ws.Cells.Item(13, 13).Value = 0.6698845136037477
$ws.Cells.Item(14, 3).Value = 0.06817264359469988
$ws.Cells.Item(14, 4).Value = 0.1554133302808509
$ws.Cells.Item(14, 5).Value = 0.07924357111372871
$ws.Cells.Item(14, 6).Value = 2.937492591539211
$ws.Cells.Item(14, 7).Value = 0.002555511581792991
$ws.Cells.Item(14, 9).Value = 2.649405013124351
$ws.Cells.Item(14, 11).Value = 2.668740061683422
$ws.Cells.Item(14, 12).Value = 0.1333260315286751
$ws.Cells.Item(14, 13).Value = 0.6647097386236709
$ws.Cells.Item(15, 3).Value = 0.06811801875181089
$ws.Cells.Item(15, 4).Value = 0.1555361793711612
$ws.Cells.Item(15, 5).Value = 0.07920718782487413
$ws.Cells.Item(15, 6).Value = 2.932026404648212
$ws.Cells.Item(15, 7).Value = 0.00255611916686535
$ws.Cells.Item(15, 9).Value = 2.64544302789119
$ws.Cells.Item(15, 11).Value = 2.651891509237885
$ws.Cells.Item(15, 12).Value = 0.1332841915720664
$ws.Cells.Item(15, 13).Value = 0.6615468476240949
$ws.Cells.Item(16, 3).Value = 0.06781190256315739
$ws.Cells.Item(16, 4).Value = 0.15625366811771
$ws.Cells.Item(16, 5).Value = 0.07900793999117184
$ws.Cells.Item(16, 6).Value = 2.901200564834852
$ws.Cells.Item(16, 7).Value = 0.00255965313251459
$ws.Cells.Item(16, 9).Value = 2.623170164728151
$ws.Cells.Item(16, 11).Value = 2.555734165151421
$ws.Cells.Item(16, 12).Value = 0.1330608842411038
$ws.Cells.Item(16, 13).Value = 0.6435257773847383
$ws.Cells.Item(17, 3).Value = 0.06763030525721092
$ws.Cells.Item(17, 4).Value = 0.1567058264279133
$ws.Cells.Item(17, 5).Value = 0.07889395947646349
$ws.Cells.Item(17, 6).Value = 2.882732849766199
$ws.Cells.Item(17, 7).Value = 0.002561867753602118
$ws.Cells.Item(17, 9).Value = 2.609890074634379
$ws.Cells.Item(17, 11).Value = 2.497096170694931
$ws.Cells.Item(17, 12).Value = 0.1329385540250954
$ws.Cells.Item(17, 13).Value = 0.6325632146227989
$ws.Cells.Item(18, 3).Value = 0.06752814121210804
$ws.Cells.Item(18, 4).Value = 0.1569702927251164
$ws.Cells.Item(18, 5).Value = 0.07883143240499635
$ws.Cells.Item(18, 6).Value = 2.872272684997938
$ws.Cells.Item(18, 7).Value = 0.002563158736059057
$ws.Cells.Item(18, 9).Value = 2.602392054515022
$ws.Cells.Item(18, 11).Value = 2.463497923630598
$ws.Cells.Item(18, 12).Value = 0.132873570590391
$ws.Cells.Item(18, 13).Value = 0.6262918647116251
$ws.Cells.Item(19, 3).Value = 0.06749394371622941
$ws.Cells.Item(19, 4).Value = 0.1570605903023612
$ws.Cells.Item(19, 5).Value = 0.07881078239361372
$ws.Cells.Item(19, 6).Value = 2.868758796706743
$ws.Cells.Item(19, 7).Value = 0.002563598797947817
$ws.Cells.Item(19, 9).Value = 2.599877396945018
$ws.Cells.Item(19, 11).Value = 2.452144217623641
$ws.Cells.Item(19, 12).Value = 0.1328524909527573
$ws.Cells.Item(19, 13).Value = 0.6241743309234096
$ws.Cells.Item(20, 3).Value = 0.06764940022639365
$ws.Cells.Item(20, 4).Value = 0.1566572380382816
$ws.Cells.Item(20, 5).Value = 0.07890577912747787
$ws.Cells.Item(20, 6).Value = 2.884681987572918
$ws.Cells.Item(20, 7).Value = 0.002561630225186205
$ws.Cells.Item(20, 9).Value = 2.611289224405084
$ws.Cells.Item(20, 11).Value = 2.503324949231285
$ws.Cells.Item(20, 12).Value = 0.132951019483599
$ws.Cells.Item(20, 13).Value = 0.6337266750018529
$ws.Cells.Item(21, 3).Value = 0.06819895971860745
$ws.Cells.Item(21, 4).Value = 0.1553546599364459
$ws.Cells.Item(21, 5).Value = 0.079261181552301
$ws.Cells.Item(21, 6).Value = 2.940122626402939
$ws.Cells.Item(21, 7).Value = 0.0025552211487596
$ws.Cells.Item(21, 9).Value = 2.65131255123886
$ws.Cells.Item(21, 11).Value = 2.676826399555921
$ws.Cells.Item(21, 12).Value = 0.1333463864737183
$ws.Cells.Item(21, 13).Value = 0.6662282762601848
$ws.Cells.Item(22, 3).Value = 0.06857828932101029
$ws.Cells.Item(22, 4).Value = 0.1545428010252365
$ws.Cells.Item(22, 5).Value = 0.07952047332899426
$ws.Cells.Item(22, 6).Value = 2.977815772899959
$ws.Cells.Item(22, 7).Value = 0.002551184507661796
$ws.Cells.Item(22, 9).Value = 2.678733256179356
$ws.Cells.Item(22, 11).Value = 2.791374387937083
$ws.Cells.Item(22, 12).Value = 0.1336528918943927
$ws.Cells.Item(22, 13).Value = 0.6877744743382408
$ws.Cells.Item(23, 3).Value = 0.06837399859168869
$ws.Cells.Item(23, 4).Value = 0.1549724810707538
$ws.Cells.Item(23, 5).Value = 0.07937961141050565
$ws.Cells.Item(23, 6).Value = 2.957563855382233
$ws.Cells.Item(23, 7).Value = 0.0025533250748877
$ws.Cells.Item(23, 9).Value = 2.663982075764594
$ws.Cells.Item(23, 11).Value = 2.730131908085298
$ws.Cells.Item(23, 12).Value = 0.1334848911967512
$ws.Cells.Item(23, 13).Value = 0.6762469226128047
$ws.Cells.Item(24, 3).Value = 0.06764076041481104
$ws.Cells.Item(24, 4).Value = 0.1566791907915928
$ws.Cells.Item(24, 5).Value = 0.07890042611336412
$ws.Cells.Item(24, 6).Value = 2.883800293007909
$ws.Cells.Item(24, 7).Value = 0.002561737556450616
$ws.Cells.Item(24, 9).Value = 2.61065624258994
$ws.Cells.Item(24, 11).Value = 2.500508564961194
$ws.Cells.Item(24, 12).Value = 0.1329453672013443
$ws.Cells.Item(24, 13).Value = 0.6332005774825262
$ws.Cells.Item(25, 3).Value = 0.06693735701151127
$ws.Cells.Item(25, 4).Value = 0.1586880835063624
$ws.Cells.Item(25, 5).Value = 0.07849943509397228
$ws.Cells.Item(25, 6).Value = 2.810416687391282
$ws.Cells.Item(25, 7).Value = 0.002571469668492707
$ws.Cells.Item(25, 9).Value = 2.558492911297563
$ws.Cells.Item(25, 11).Value = 2.257740859083697
$ws.Cells.Item(25, 12).Value = 0.1325688690902922
$ws.Cells.Item(25, 13).Value = 0.588068042496019
